$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ETLE")
$ws.Range("B2").Value = -5
$wb.Worksheets.Item("About").Activate()
